$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before the old column C ("TIPO TRAMITE"), shifting
# everything from C onward two places to the right (old M -> new O).
$ws.Columns("C:D").Insert()

# Fill the two new header cells. Set D1 before C1 so the shared-string table
# picks up "BL HOUSE" (index 13) ahead of "BL MASTER" (index 14), matching
# how the original edit was authored.
$ws.Range("D1").Value = "BL HOUSE"
$ws.Range("C1").Value = "BL MASTER"

# Give the two new columns an explicit (non bestFit) width matching column B.
$ws.Range("C1:D1").ColumnWidth = 15.5

# Move the selection to A4.
$ws.Range("A4").Select()

# The sheet-level AutoFilter range needs to grow from A1:M1 to A1:O1. Toggling
# AutoFilterMode off first avoids Range.AutoFilter() acting as an OFF-toggle.
$ws.AutoFilterMode = $false
$ws.Range("A1:O1").AutoFilter()

# The workbook-level _FilterDatabase defined name also needs to track the
# wider range.
$wb.Names.Item("Hoja1!_FilterDatabase").RefersTo = "=Hoja1!`$A`$1:`$O`$1"
